{"js": "// Apply the \"Added many more features\" edit to the Don Juan's Peppers\n// review document using the Word JavaScript API (Office.js).\n//\n// Strategy: each change is a self-contained text substitution of a whole\n// run's text. We locate the old text with Range.search (exact, whole-word\n// off so punctuation-containing strings still match) and replace it via\n// insertText(..., Word.InsertLocation.replace), which swaps the text of\n// the found range while leaving its run formatting (bold/italic/etc.)\n// untouched.\n\nasync function replaceAll(context, oldText, newText) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// 1. Page heading (also reused verbatim further down in a bold run).\nawait replaceAll(\n  context,\n  \"Play Don Juan's Peppers free online slot game\",\n  \"Play Don Juan's Peppers for Free\"\n);\n\n// 2. \"What we like\" bullet list.\nawait replaceAll(\n  context,\n  \"Immersive, Mexican-themed gameplay\",\n  \"Immersive Mexican theme\"\n);\nawait replaceAll(\n  context,\n  \"Big win opportunities with scatter and wild symbols\",\n  \"Artfully designed symbols\"\n);\nawait replaceAll(\n  context,\n  \"Free spins mode with extended chances\",\n  \"Generous free spins mode\"\n);\nawait replaceAll(\n  context,\n  \"Gamble feature allows players to double their winnings\",\n  \"Big win opportunities\"\n);\n\n// 3. \"What we don't like\" bullet list.\nawait replaceAll(\n  context,\n  \"Limited amount of bonus features\",\n  \"Limited bonus features\"\n);\nawait replaceAll(\n  context,\n  \"Graphic design might not suit everyone's taste\",\n  \"No progressive jackpot\"\n);\n\n// 4. Closing meta description (italic run).\nawait replaceAll(\n  context,\n  \"Read our Don Juan's Peppers review and play for free. Enjoy a Mexican-themed slot game with big win opportunities and a generous free spins mode.\",\n  \"Read our review of Don Juan's Peppers and play for free. Experience the Mexican-themed slot game with big win opportunities.\"\n);\n", "ps1": "# Apply the \"Added many more features\" edit to the Don Juan's Peppers\n# review document using the Word COM object model.\n#\n# We use Find.Execute() only to *locate* each target run's text, then\n# assign straight to Range.Text to swap its contents. Doing the\n# replacement via direct Range.Text assignment (rather than\n# Find.Execute(..., Replace:=wdReplaceAll)) sidesteps Word's smart-quote\n# autocorrect, which would otherwise mangle the straight apostrophes\n# that appear in both the old and new copy.\n\nfunction Replace-AllText($doc, $oldText, $newText) {\n    $count = 0\n    while ($true) {\n        $rng = $doc.Content\n        $rng.Find.Text = $oldText\n        $rng.Find.Forward = $true\n        $rng.Find.Wrap = 0\n        $found = $rng.Find.Execute()\n        if (-not $found) { break }\n        $rng.Text = $newText\n        $count = $count + 1\n    }\n    return $count\n}\n\n$d = $word.ActiveDocument\n\n# 1. Page heading (also reused verbatim further down in a bold run).\nReplace-AllText $d \"Play Don Juan's Peppers free online slot game\" \"Play Don Juan's Peppers for Free\" | Out-Null\n\n# 2. \"What we like\" bullet list.\nReplace-AllText $d \"Immersive, Mexican-themed gameplay\" \"Immersive Mexican theme\" | Out-Null\nReplace-AllText $d \"Big win opportunities with scatter and wild symbols\" \"Artfully designed symbols\" | Out-Null\nReplace-AllText $d \"Free spins mode with extended chances\" \"Generous free spins mode\" | Out-Null\nReplace-AllText $d \"Gamble feature allows players to double their winnings\" \"Big win opportunities\" | Out-Null\n\n# 3. \"What we don't like\" bullet list.\nReplace-AllText $d \"Limited amount of bonus features\" \"Limited bonus features\" | Out-Null\nReplace-AllText $d \"Graphic design might not suit everyone's taste\" \"No progressive jackpot\" | Out-Null\n\n# 4. Closing meta description (italic run).\nReplace-AllText $d \"Read our Don Juan's Peppers review and play for free. Enjoy a Mexican-themed slot game with big win opportunities and a generous free spins mode.\" \"Read our review of Don Juan's Peppers and play for free. Experience the Mexican-themed slot game with big win opportunities.\" | Out-Null\n"}
